{"js": "// Add a new \"Universal Oled Lib\" entry (title + link) at the very top of the\n// document body. The \"_GoBack\" bookmark that used to sit in the last\n// (now-empty) paragraph moves to the end of the new title paragraph.\n\nconst body = context.document.body;\n\n// Insert the two new paragraphs at the start of the body. Each call targets\n// Word.InsertLocation.start, so insert the second-desired-line first and the\n// first-desired-line last, leaving them in the correct final order:\n//   1. \"Universal Oled Lib \"\n//   2. \"https://code.google.com/archive/p/u8glib/\"\nbody.insertParagraph(\"https://code.google.com/archive/p/u8glib/\", Word.InsertLocation.start);\nbody.insertParagraph(\"Universal Oled Lib \", Word.InsertLocation.start);\nawait context.sync();\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The \"Universal Oled Lib \" paragraph is now the first paragraph in the body.\nconst titleParagraph = paragraphs.items[0];\nconst titleEnd = titleParagraph.getRange(Word.RangeLocation.end);\n\n// Relocate the \"_GoBack\" bookmark: delete the old one (it was on the trailing\n// empty paragraph) and recreate it at the end of the title paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\ntitleEnd.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Add a new \"Universal Oled Lib\" entry (title + link) at the very top of the\n# document. The \"_GoBack\" bookmark that used to live on the last (now-empty)\n# paragraph moves to the end of the new title paragraph.\n\n$d = $word.ActiveDocument\n\n# Insert the two new lines at the very start of the document. A placeholder\n# character (\"X\") is appended right after the title text so that the bookmark\n# below can be anchored at a safe (non-paragraph-edge) text position; it is\n# removed again immediately afterwards.\n$startRange = $d.Range(0, 0)\n$startRange.InsertBefore(\"Universal Oled Lib X\" + [char]13 + \"https://code.google.com/archive/p/u8glib/\" + [char]13)\n\n# Drop the old \"_GoBack\" bookmark (it lived on the trailing empty paragraph).\n$oldBookmark = $d.Bookmarks.Item(\"_GoBack\")\n$oldBookmark.Delete()\n\n# Re-create \"_GoBack\" collapsed right before the placeholder \"X\", i.e. right\n# after \"Universal Oled Lib \".\n$titleLen = \"Universal Oled Lib \".Length\n$bookmarkRange = $d.Range($titleLen, $titleLen)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange) | Out-Null\n\n# Remove the placeholder character now that the bookmark is anchored.\n$placeholderRange = $d.Range($titleLen, $titleLen + 1)\n$placeholderRange.Delete()\n"}
